$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value2 = 1242.25
$ws.Cells.Item(2, 9).Value2 = 388
$ws.Cells.Item(2, 11).Value2 = 388
$ws.Cells.Item(2, 13).Value2 = -275
$ws.Cells.Item(17, 8).Value2 = 112243.664
$ws.Cells.Item(17, 10).Value2 = 126136.625
$ws.Cells.Item(17, 12).Value2 = 378409.875
$ws.Cells.Item(17, 14).Value2 = -378745.875
$ws.Cells.Item(38, 8).Value2 = 1680.1538
$ws.Cells.Item(38, 9).Value2 = 249.4
$ws.Cells.Item(38, 10).Value2 = 6449.3335
$ws.Cells.Item(38, 11).Value2 = 748.2
$ws.Cells.Item(38, 12).Value2 = 19348.0005
$ws.Cells.Item(38, 13).Value2 = -376.2
$ws.Cells.Item(38, 14).Value2 = -20092.0005
$ws.Cells.Item(40, 8).Value2 = 71430500
$ws.Cells.Item(40, 9).Value2 = 1449.5
$ws.Cells.Item(40, 10).Value2 = 100002110
$ws.Cells.Item(40, 11).Value2 = 1449.5
$ws.Cells.Item(40, 12).Value2 = 100002110
$ws.Cells.Item(40, 13).Value2 = -1274.5
$ws.Cells.Item(40, 14).Value2 = -100002460
$ws.Cells.Item(75, 8).Value2 = 64749.125
$ws.Cells.Item(75, 9).Value2 = 39999
$ws.Cells.Item(75, 10).Value2 = 68284.86
$ws.Cells.Item(75, 11).Value2 = 39999
$ws.Cells.Item(75, 12).Value2 = 68284.86
$ws.Cells.Item(75, 13).Value2 = -39063
$ws.Cells.Item(75, 14).Value2 = -70156.86
$ws.Cells.Item(78, 8).Value2 = 64749.125
$ws.Cells.Item(78, 9).Value2 = 39999
$ws.Cells.Item(78, 10).Value2 = 68284.86
$ws.Cells.Item(78, 11).Value2 = 119997
$ws.Cells.Item(78, 12).Value2 = 204854.58
$ws.Cells.Item(78, 13).Value2 = -115317
$ws.Cells.Item(78, 14).Value2 = -214214.58
$ws.Cells.Item(95, 8).Value2 = 53539.332
$ws.Cells.Item(95, 10).Value2 = 53539.332
$ws.Cells.Item(95, 12).Value2 = 53539.332
$ws.Cells.Item(95, 14).Value2 = -59031.332
$ws.Cells.Item(103, 8).Value2 = 35716070
$ws.Cells.Item(103, 9).Value2 = 1036.6
$ws.Cells.Item(103, 10).Value2 = 55557760
$ws.Cells.Item(103, 11).Value2 = 3109.8
$ws.Cells.Item(103, 12).Value2 = 166673280
$ws.Cells.Item(103, 13).Value2 = -2523.8
$ws.Cells.Item(103, 14).Value2 = -166674452
$ws.Cells.Item(105, 8).Value2 = 52981.5
$ws.Cells.Item(105, 10).Value2 = 52981.5
$ws.Cells.Item(105, 12).Value2 = 52981.5
$ws.Cells.Item(105, 14).Value2 = -59969.5
$ws.Cells.Item(109, 8).Value2 = 0
$ws.Cells.Item(109, 10).Value2 = 0
$ws.Cells.Item(109, 12).Value2 = 0
$ws.Cells.Item(109, 14).Value2 = $null
$ws.Cells.Item(132, 8).Value2 = 2414
$ws.Cells.Item(132, 9).Value2 = 2427.516
$ws.Cells.Item(132, 10).Value2 = 1995
$ws.Cells.Item(132, 11).Value2 = 7282.548000000001
$ws.Cells.Item(132, 12).Value2 = 5985
$ws.Cells.Item(132, 13).Value2 = -4752.548000000001
$ws.Cells.Item(132, 14).Value2 = -11045
$ws.Cells.Item(133, 8).Value2 = 118636.5
$ws.Cells.Item(133, 10).Value2 = 118636.5
$ws.Cells.Item(133, 12).Value2 = 118636.5
$ws.Cells.Item(133, 14).Value2 = -128756.5
$ws.Cells.Item(135, 8).Value2 = 3860.2727
$ws.Cells.Item(135, 9).Value2 = 2380.4285
$ws.Cells.Item(135, 10).Value2 = 6450
$ws.Cells.Item(135, 11).Value2 = 21423.8565
$ws.Cells.Item(135, 12).Value2 = 58050
$ws.Cells.Item(135, 13).Value2 = -18888.8565
$ws.Cells.Item(135, 14).Value2 = -63120
$ws.Cells.Item(138, 8).Value2 = 4718.722
$ws.Cells.Item(138, 9).Value2 = 4450.3335
$ws.Cells.Item(138, 10).Value2 = 4772.4
$ws.Cells.Item(138, 11).Value2 = 13351.0005
$ws.Cells.Item(138, 12).Value2 = 14317.2
$ws.Cells.Item(138, 13).Value2 = -8211.000499999998
$ws.Cells.Item(138, 14).Value2 = -24597.2
$ws.Cells.Item(141, 8).Value2 = 4267.0835
$ws.Cells.Item(141, 9).Value2 = 4880
$ws.Cells.Item(141, 10).Value2 = 1202.5
$ws.Cells.Item(141, 11).Value2 = 14640
$ws.Cells.Item(141, 12).Value2 = 3607.5
$ws.Cells.Item(141, 13).Value2 = -9460
$ws.Cells.Item(141, 14).Value2 = -13967.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value2 = 1017.0357
$ws.Cells.Item(2, 9).Value2 = 976.8421
$ws.Cells.Item(2, 10).Value2 = 1101.8889
$ws.Cells.Item(2, 11).Value2 = 976.8421
$ws.Cells.Item(2, 12).Value2 = 1101.8889
$ws.Cells.Item(2, 13).Value2 = -863.8421
$ws.Cells.Item(2, 14).Value2 = -1327.8889
$ws.Cells.Item(32, 8).Value2 = 9433.540999999999
$ws.Cells.Item(32, 9).Value2 = 9273.062
$ws.Cells.Item(32, 11).Value2 = 9273.062
$ws.Cells.Item(32, 13).Value2 = -8986.062
$ws.Cells.Item(36, 8).Value2 = 24519.8
$ws.Cells.Item(36, 9).Value2 = 10200
$ws.Cells.Item(36, 10).Value2 = 45999.5
$ws.Cells.Item(36, 11).Value2 = 10200
$ws.Cells.Item(36, 12).Value2 = 45999.5
$ws.Cells.Item(36, 13).Value2 = -9854
$ws.Cells.Item(36, 14).Value2 = -46691.5
$ws.Cells.Item(45, 8).Value2 = 2557.9412
$ws.Cells.Item(45, 9).Value2 = 1114.2307
$ws.Cells.Item(45, 11).Value2 = 1114.2307
$ws.Cells.Item(45, 13).Value2 = -737.2307000000001
$ws.Cells.Item(58, 8).Value2 = 49999.5
$ws.Cells.Item(58, 10).Value2 = 49999.5
$ws.Cells.Item(58, 12).Value2 = 49999.5
$ws.Cells.Item(58, 14).Value2 = -50859.5
$ws.Cells.Item(61, 8).Value2 = 6567644.5
$ws.Cells.Item(61, 9).Value2 = 9096743
$ws.Cells.Item(61, 11).Value2 = 9096743
$ws.Cells.Item(61, 13).Value2 = -9096531
$ws.Cells.Item(74, 8).Value2 = 2552.5789
$ws.Cells.Item(74, 9).Value2 = 2527.7222
$ws.Cells.Item(74, 10).Value2 = 3000
$ws.Cells.Item(74, 11).Value2 = 2527.7222
$ws.Cells.Item(74, 12).Value2 = 3000
$ws.Cells.Item(74, 13).Value2 = -1653.7222
$ws.Cells.Item(74, 14).Value2 = -4748
$ws.Cells.Item(77, 8).Value2 = 2552.5789
$ws.Cells.Item(77, 9).Value2 = 2527.7222
$ws.Cells.Item(77, 10).Value2 = 3000
$ws.Cells.Item(77, 11).Value2 = 12638.611
$ws.Cells.Item(77, 12).Value2 = 15000
$ws.Cells.Item(77, 13).Value2 = -8270.611000000001
$ws.Cells.Item(77, 14).Value2 = -23736
$ws.Cells.Item(102, 8).Value2 = 2598.652
$ws.Cells.Item(102, 9).Value2 = 2012.8572
$ws.Cells.Item(102, 11).Value2 = 2012.8572
$ws.Cells.Item(102, 13).Value2 = -390.8571999999999
$ws.Cells.Item(110, 8).Value2 = 4212.7144
$ws.Cells.Item(110, 9).Value2 = 3946.3928
$ws.Cells.Item(110, 10).Value2 = 5278
$ws.Cells.Item(110, 11).Value2 = 3946.3928
$ws.Cells.Item(110, 12).Value2 = 5278
$ws.Cells.Item(110, 13).Value2 = -1901.3928
$ws.Cells.Item(110, 14).Value2 = -9368
$ws.Cells.Item(116, 8).Value2 = 1017.0357
$ws.Cells.Item(116, 9).Value2 = 976.8421
$ws.Cells.Item(116, 10).Value2 = 1101.8889
$ws.Cells.Item(116, 11).Value2 = 976.8421
$ws.Cells.Item(116, 12).Value2 = 1101.8889
$ws.Cells.Item(116, 13).Value2 = 1317.1579
$ws.Cells.Item(116, 14).Value2 = -5689.8889
$ws.Cells.Item(122, 8).Value2 = 4552.7666
$ws.Cells.Item(122, 9).Value2 = 3583.52
$ws.Cells.Item(122, 11).Value2 = 10750.56
$ws.Cells.Item(122, 13).Value2 = -8300.559999999999
$ws.Cells.Item(132, 8).Value2 = 2780574.8
$ws.Cells.Item(132, 9).Value2 = 2828.375
$ws.Cells.Item(132, 11).Value2 = 8485.125
$ws.Cells.Item(132, 13).Value2 = -5955.125
$ws.Cells.Item(136, 8).Value2 = 6567644.5
$ws.Cells.Item(136, 9).Value2 = 9096743
$ws.Cells.Item(136, 11).Value2 = 27290229
$ws.Cells.Item(136, 13).Value2 = -27287679
$ws.Cells.Item(137, 8).Value2 = 109093.8
$ws.Cells.Item(137, 10).Value2 = 109093.8
$ws.Cells.Item(137, 12).Value2 = 109093.8
$ws.Cells.Item(137, 14).Value2 = -119293.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value2 = 1017.0357
$ws.Cells.Item(3, 9).Value2 = 976.8421
$ws.Cells.Item(3, 10).Value2 = 1101.8889
$ws.Cells.Item(3, 11).Value2 = 976.8421
$ws.Cells.Item(3, 12).Value2 = 1101.8889
$ws.Cells.Item(3, 13).Value2 = -862.8421
$ws.Cells.Item(3, 14).Value2 = -1329.8889
$ws.Cells.Item(86, 8).Value2 = 6849.4
$ws.Cells.Item(86, 9).Value2 = 3998.75
$ws.Cells.Item(86, 10).Value2 = 8749.833000000001
$ws.Cells.Item(86, 11).Value2 = 3998.75
$ws.Cells.Item(86, 12).Value2 = 8749.833000000001
$ws.Cells.Item(86, 13).Value2 = -2875.75
$ws.Cells.Item(86, 14).Value2 = -10995.833
$ws.Cells.Item(89, 8).Value2 = 6849.4
$ws.Cells.Item(89, 9).Value2 = 3998.75
$ws.Cells.Item(89, 10).Value2 = 8749.833000000001
$ws.Cells.Item(89, 11).Value2 = 19993.75
$ws.Cells.Item(89, 12).Value2 = 43749.165
$ws.Cells.Item(89, 13).Value2 = -14377.75
$ws.Cells.Item(89, 14).Value2 = -54981.165
$ws.Cells.Item(94, 8).Value2 = 2442.8215
$ws.Cells.Item(94, 9).Value2 = 2937.2104
$ws.Cells.Item(94, 10).Value2 = 1399.1111
$ws.Cells.Item(94, 11).Value2 = 2937.2104
$ws.Cells.Item(94, 12).Value2 = 1399.1111
$ws.Cells.Item(94, 13).Value2 = -2486.2104
$ws.Cells.Item(94, 14).Value2 = -2301.1111
$ws.Cells.Item(99, 8).Value2 = 2642.111
$ws.Cells.Item(99, 9).Value2 = 1000
$ws.Cells.Item(99, 10).Value2 = 2847.375
$ws.Cells.Item(99, 11).Value2 = 1000
$ws.Cells.Item(99, 12).Value2 = 2847.375
$ws.Cells.Item(99, 13).Value2 = 498
$ws.Cells.Item(99, 14).Value2 = -5843.375
$ws.Cells.Item(105, 8).Value2 = 808487.3
$ws.Cells.Item(105, 9).Value2 = 1430800.1
$ws.Cells.Item(105, 10).Value2 = 8370.857
$ws.Cells.Item(105, 11).Value2 = 1430800.1
$ws.Cells.Item(105, 12).Value2 = 8370.857
$ws.Cells.Item(105, 13).Value2 = -1429053.1
$ws.Cells.Item(105, 14).Value2 = -11864.857
$ws.Cells.Item(107, 8).Value2 = 7512.222
$ws.Cells.Item(107, 9).Value2 = 7638.75
$ws.Cells.Item(107, 10).Value2 = 6500
$ws.Cells.Item(107, 11).Value2 = 7638.75
$ws.Cells.Item(107, 12).Value2 = 6500
$ws.Cells.Item(107, 13).Value2 = -5718.75
$ws.Cells.Item(107, 14).Value2 = -10340
$ws.Cells.Item(132, 8).Value2 = 112311
$ws.Cells.Item(132, 10).Value2 = 112311
$ws.Cells.Item(132, 12).Value2 = 112311
$ws.Cells.Item(132, 14).Value2 = -122431
$ws.Cells.Item(140, 8).Value2 = 139985.75
$ws.Cells.Item(140, 10).Value2 = 139985.75
$ws.Cells.Item(140, 12).Value2 = 139985.75
$ws.Cells.Item(140, 14).Value2 = -150345.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 17859736
$ws.Cells.Item(31, 9).Value2 = 18520976
$ws.Cells.Item(31, 10).Value2 = 6209
$ws.Cells.Item(31, 11).Value2 = 18520976
$ws.Cells.Item(31, 12).Value2 = 6209
$ws.Cells.Item(31, 13).Value2 = -18520681
$ws.Cells.Item(31, 14).Value2 = -6799
$ws.Cells.Item(34, 8).Value2 = 17859736
$ws.Cells.Item(34, 9).Value2 = 18520976
$ws.Cells.Item(34, 10).Value2 = 6209
$ws.Cells.Item(34, 11).Value2 = 18520976
$ws.Cells.Item(34, 12).Value2 = 6209
$ws.Cells.Item(34, 13).Value2 = -18520774
$ws.Cells.Item(34, 14).Value2 = -6613
$ws.Cells.Item(99, 8).Value2 = 12865.441
$ws.Cells.Item(99, 9).Value2 = 6939.524
$ws.Cells.Item(99, 10).Value2 = 22438.076
$ws.Cells.Item(99, 11).Value2 = 6939.524
$ws.Cells.Item(99, 12).Value2 = 22438.076
$ws.Cells.Item(99, 13).Value2 = -5441.524
$ws.Cells.Item(99, 14).Value2 = -25434.076
$ws.Cells.Item(107, 8).Value2 = 1982.1538
$ws.Cells.Item(107, 9).Value2 = 623.5
$ws.Cells.Item(107, 10).Value2 = 4156
$ws.Cells.Item(107, 11).Value2 = 623.5
$ws.Cells.Item(107, 12).Value2 = 4156
$ws.Cells.Item(107, 13).Value2 = 1296.5
$ws.Cells.Item(107, 14).Value2 = -7996
$ws.Cells.Item(120, 8).Value2 = 59038
$ws.Cells.Item(120, 10).Value2 = 59038
$ws.Cells.Item(120, 12).Value2 = 59038
$ws.Cells.Item(120, 14).Value2 = -66296
$ws.Cells.Item(126, 8).Value2 = 12865.441
$ws.Cells.Item(126, 9).Value2 = 6939.524
$ws.Cells.Item(126, 10).Value2 = 22438.076
$ws.Cells.Item(126, 11).Value2 = 20818.572
$ws.Cells.Item(126, 12).Value2 = 67314.228
$ws.Cells.Item(126, 13).Value2 = -18348.572
$ws.Cells.Item(126, 14).Value2 = -72254.228
$ws.Cells.Item(132, 8).Value2 = 2081.4583
$ws.Cells.Item(132, 9).Value2 = 2010.762
$ws.Cells.Item(132, 10).Value2 = 2576.3333
$ws.Cells.Item(132, 11).Value2 = 6032.286
$ws.Cells.Item(132, 12).Value2 = 7728.999899999999
$ws.Cells.Item(132, 13).Value2 = -3502.286
$ws.Cells.Item(132, 14).Value2 = -12788.9999
$ws.Cells.Item(134, 8).Value2 = 306
$ws.Cells.Item(134, 9).Value2 = 306
$ws.Cells.Item(134, 11).Value2 = 918
$ws.Cells.Item(134, 13).Value2 = 1617
$ws.Cells.Item(135, 8).Value2 = 71955
$ws.Cells.Item(135, 9).Value2 = 71955
$ws.Cells.Item(135, 11).Value2 = 71955
$ws.Cells.Item(135, 13).Value2 = -66885
$ws.Cells.Item(141, 8).Value2 = 327575
$ws.Cells.Item(141, 10).Value2 = 365370.38
$ws.Cells.Item(141, 12).Value2 = 365370.38
$ws.Cells.Item(141, 14).Value2 = -375730.38

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value2 = 5070.533
$ws.Cells.Item(12, 10).Value2 = 4881
$ws.Cells.Item(12, 12).Value2 = 14643
$ws.Cells.Item(12, 14).Value2 = -14989
$ws.Cells.Item(69, 8).Value2 = 10557
$ws.Cells.Item(69, 9).Value2 = 3625.4285
$ws.Cells.Item(69, 10).Value2 = 16622.125
$ws.Cells.Item(69, 11).Value2 = 10876.2855
$ws.Cells.Item(69, 12).Value2 = 49866.375
$ws.Cells.Item(69, 13).Value2 = -10065.2855
$ws.Cells.Item(69, 14).Value2 = -51488.375
$ws.Cells.Item(72, 8).Value2 = 10557
$ws.Cells.Item(72, 9).Value2 = 3625.4285
$ws.Cells.Item(72, 10).Value2 = 16622.125
$ws.Cells.Item(72, 11).Value2 = 32628.8565
$ws.Cells.Item(72, 12).Value2 = 149599.125
$ws.Cells.Item(72, 13).Value2 = -28572.8565
$ws.Cells.Item(72, 14).Value2 = -157711.125
$ws.Cells.Item(92, 8).Value2 = 284.91666
$ws.Cells.Item(92, 10).Value2 = 322
$ws.Cells.Item(92, 12).Value2 = 966
$ws.Cells.Item(92, 14).Value2 = -3462
$ws.Cells.Item(129, 8).Value2 = 3773.842
$ws.Cells.Item(129, 9).Value2 = 3325.5833
$ws.Cells.Item(129, 10).Value2 = 4542.2856
$ws.Cells.Item(129, 11).Value2 = 9976.749899999999
$ws.Cells.Item(129, 12).Value2 = 13626.8568
$ws.Cells.Item(129, 13).Value2 = -4976.749899999999
$ws.Cells.Item(129, 14).Value2 = -23626.8568
$ws.Cells.Item(132, 8).Value2 = 4700.533
$ws.Cells.Item(132, 10).Value2 = 4732.3335
$ws.Cells.Item(132, 12).Value2 = 42591.0015
$ws.Cells.Item(132, 14).Value2 = -47651.0015
$ws.Cells.Item(134, 8).Value2 = 10241.917
$ws.Cells.Item(134, 9).Value2 = 3655.889
$ws.Cells.Item(134, 11).Value2 = 10967.667
$ws.Cells.Item(134, 13).Value2 = -5897.667000000001
$ws.Cells.Item(137, 8).Value2 = 8878.468999999999
$ws.Cells.Item(137, 9).Value2 = 5350
$ws.Cells.Item(137, 10).Value2 = 10726.714
$ws.Cells.Item(137, 11).Value2 = 16050
$ws.Cells.Item(137, 12).Value2 = 32180.142
$ws.Cells.Item(137, 13).Value2 = -10950
$ws.Cells.Item(137, 14).Value2 = -42380.142
$ws.Cells.Item(140, 8).Value2 = 4409.6875
$ws.Cells.Item(140, 9).Value2 = 1830.4445
$ws.Cells.Item(140, 11).Value2 = 5491.333500000001
$ws.Cells.Item(140, 13).Value2 = -311.3335000000006
$ws.Cells.Item(141, 8).Value2 = 3921.5
$ws.Cells.Item(141, 9).Value2 = 3705.8
$ws.Cells.Item(141, 10).Value2 = 5000
$ws.Cells.Item(141, 11).Value2 = 11117.4
$ws.Cells.Item(141, 12).Value2 = 15000
$ws.Cells.Item(141, 13).Value2 = -5937.400000000001
$ws.Cells.Item(141, 14).Value2 = -25360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value2 = 7735
$ws.Cells.Item(80, 9).Value2 = 4701.6665
$ws.Cells.Item(80, 11).Value2 = 4701.6665
$ws.Cells.Item(80, 13).Value2 = -3703.6665
$ws.Cells.Item(83, 8).Value2 = 7735
$ws.Cells.Item(83, 9).Value2 = 4701.6665
$ws.Cells.Item(83, 11).Value2 = 23508.3325
$ws.Cells.Item(83, 13).Value2 = -18516.3325
$ws.Cells.Item(97, 8).Value2 = 1245.0312
$ws.Cells.Item(97, 9).Value2 = 1147.44
$ws.Cells.Item(97, 10).Value2 = 1593.5714
$ws.Cells.Item(97, 11).Value2 = 1147.44
$ws.Cells.Item(97, 12).Value2 = 1593.5714
$ws.Cells.Item(97, 13).Value2 = -651.4400000000001
$ws.Cells.Item(97, 14).Value2 = -2585.5714
$ws.Cells.Item(102, 8).Value2 = 3018.3076
$ws.Cells.Item(102, 9).Value2 = 2104.75
$ws.Cells.Item(102, 10).Value2 = 4480
$ws.Cells.Item(102, 11).Value2 = 2104.75
$ws.Cells.Item(102, 12).Value2 = 4480
$ws.Cells.Item(102, 13).Value2 = -482.75
$ws.Cells.Item(102, 14).Value2 = -7724
$ws.Cells.Item(122, 8).Value2 = 3461.6667
$ws.Cells.Item(122, 9).Value2 = 3540.6667
$ws.Cells.Item(122, 10).Value2 = 3277.3333
$ws.Cells.Item(122, 11).Value2 = 10622.0001
$ws.Cells.Item(122, 12).Value2 = 9831.999899999999
$ws.Cells.Item(122, 13).Value2 = -8172.000100000001
$ws.Cells.Item(122, 14).Value2 = -14731.9999
$ws.Cells.Item(132, 8).Value2 = 3336931.2
$ws.Cells.Item(132, 9).Value2 = 3651.577
$ws.Cells.Item(132, 11).Value2 = 10954.731
$ws.Cells.Item(132, 13).Value2 = -8424.731

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value2 = 1749.3334
$ws.Cells.Item(16, 9).Value2 = 1749.3334
$ws.Cells.Item(16, 10).Value2 = 0
$ws.Cells.Item(16, 11).Value2 = 1749.3334
$ws.Cells.Item(16, 12).Value2 = 0
$ws.Cells.Item(16, 13).Value2 = -1579.3334
$ws.Cells.Item(16, 14).Value2 = $null
$ws.Cells.Item(46, 8).Value2 = 971.13336
$ws.Cells.Item(46, 9).Value2 = 559.3
$ws.Cells.Item(46, 11).Value2 = 559.3
$ws.Cells.Item(46, 13).Value2 = -371.3
$ws.Cells.Item(61, 8).Value2 = 3729.9614
$ws.Cells.Item(61, 9).Value2 = 3435.6365
$ws.Cells.Item(61, 11).Value2 = 3435.6365
$ws.Cells.Item(61, 13).Value2 = -3233.6365
$ws.Cells.Item(64, 8).Value2 = 0
$ws.Cells.Item(64, 10).Value2 = 0
$ws.Cells.Item(64, 12).Value2 = 0
$ws.Cells.Item(64, 14).Value2 = $null
$ws.Cells.Item(67, 8).Value2 = 0
$ws.Cells.Item(67, 10).Value2 = 0
$ws.Cells.Item(67, 12).Value2 = 0
$ws.Cells.Item(67, 14).Value2 = $null
$ws.Cells.Item(70, 8).Value2 = 17999
$ws.Cells.Item(70, 10).Value2 = 17999
$ws.Cells.Item(70, 12).Value2 = 17999
$ws.Cells.Item(70, 14).Value2 = -18539
$ws.Cells.Item(73, 8).Value2 = 17999
$ws.Cells.Item(73, 10).Value2 = 17999
$ws.Cells.Item(73, 12).Value2 = 17999
$ws.Cells.Item(73, 14).Value2 = -19871
$ws.Cells.Item(93, 8).Value2 = 2780013.2
$ws.Cells.Item(93, 9).Value2 = 1251.7858
$ws.Cells.Item(93, 10).Value2 = 9263790
$ws.Cells.Item(93, 11).Value2 = 1251.7858
$ws.Cells.Item(93, 12).Value2 = 9263790
$ws.Cells.Item(93, 13).Value2 = -3.785800000000108
$ws.Cells.Item(93, 14).Value2 = -9266286
$ws.Cells.Item(103, 8).Value2 = 51775.5
$ws.Cells.Item(103, 10).Value2 = 51775.5
$ws.Cells.Item(103, 12).Value2 = 51775.5
$ws.Cells.Item(103, 14).Value2 = -54119.5
$ws.Cells.Item(113, 8).Value2 = 3729.9614
$ws.Cells.Item(113, 9).Value2 = 3435.6365
$ws.Cells.Item(113, 11).Value2 = 3435.6365
$ws.Cells.Item(113, 13).Value2 = -1265.6365
$ws.Cells.Item(122, 8).Value2 = 3580.1707
$ws.Cells.Item(122, 9).Value2 = 3556.8108
$ws.Cells.Item(122, 10).Value2 = 3796.25
$ws.Cells.Item(122, 11).Value2 = 10670.4324
$ws.Cells.Item(122, 12).Value2 = 11388.75
$ws.Cells.Item(122, 13).Value2 = -8220.432400000002
$ws.Cells.Item(122, 14).Value2 = -16288.75
$ws.Cells.Item(132, 8).Value2 = 4464.387
$ws.Cells.Item(132, 9).Value2 = 2829.7693
$ws.Cells.Item(132, 10).Value2 = 5644.9443
$ws.Cells.Item(132, 11).Value2 = 8489.3079
$ws.Cells.Item(132, 12).Value2 = 16934.8329
$ws.Cells.Item(132, 13).Value2 = -5959.3079
$ws.Cells.Item(132, 14).Value2 = -21994.8329
$ws.Cells.Item(133, 8).Value2 = 162499.5
$ws.Cells.Item(133, 10).Value2 = 162499.5
$ws.Cells.Item(133, 12).Value2 = 162499.5
$ws.Cells.Item(133, 14).Value2 = -167559.5
$ws.Cells.Item(136, 8).Value2 = 3810.9546
$ws.Cells.Item(136, 9).Value2 = 3717.1
$ws.Cells.Item(136, 10).Value2 = 4749.5
$ws.Cells.Item(136, 11).Value2 = 11151.3
$ws.Cells.Item(136, 12).Value2 = 14248.5
$ws.Cells.Item(136, 13).Value2 = -8601.299999999999
$ws.Cells.Item(136, 14).Value2 = -19348.5
$ws.Cells.Item(139, 8).Value2 = 75000
$ws.Cells.Item(139, 9).Value2 = 75000
$ws.Cells.Item(139, 10).Value2 = 0
$ws.Cells.Item(139, 11).Value2 = 75000
$ws.Cells.Item(139, 12).Value2 = 0
$ws.Cells.Item(139, 13).Value2 = -69860
$ws.Cells.Item(139, 14).Value2 = $null

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value2 = 2378.7693
$ws.Cells.Item(100, 9).Value2 = 1416
$ws.Cells.Item(100, 10).Value2 = 3919.2
$ws.Cells.Item(100, 11).Value2 = 2832
$ws.Cells.Item(100, 12).Value2 = 7838.4
$ws.Cells.Item(100, 13).Value2 = -2291
$ws.Cells.Item(100, 14).Value2 = -8920.4
$ws.Cells.Item(107, 8).Value2 = 4888.36
$ws.Cells.Item(107, 9).Value2 = 2088.625
$ws.Cells.Item(107, 11).Value2 = 6265.875
$ws.Cells.Item(107, 13).Value2 = -4345.875
$ws.Cells.Item(113, 8).Value2 = 1252.3334
$ws.Cells.Item(113, 9).Value2 = 1033.25
$ws.Cells.Item(113, 10).Value2 = 1361.875
$ws.Cells.Item(113, 11).Value2 = 3099.75
$ws.Cells.Item(113, 12).Value2 = 4085.625
$ws.Cells.Item(113, 13).Value2 = -929.75
$ws.Cells.Item(113, 14).Value2 = -8425.625
$ws.Cells.Item(122, 8).Value2 = 4066.4443
$ws.Cells.Item(122, 9).Value2 = 3639.8
$ws.Cells.Item(122, 10).Value2 = 4599.75
$ws.Cells.Item(122, 11).Value2 = 10919.4
$ws.Cells.Item(122, 12).Value2 = 13799.25
$ws.Cells.Item(122, 13).Value2 = -8469.400000000001
$ws.Cells.Item(122, 14).Value2 = -18699.25
$ws.Cells.Item(123, 8).Value2 = 196166.33
$ws.Cells.Item(123, 10).Value2 = 196166.33
$ws.Cells.Item(123, 12).Value2 = 196166.33
$ws.Cells.Item(123, 14).Value2 = -205966.33
$ws.Cells.Item(132, 8).Value2 = 280076.38
$ws.Cells.Item(132, 9).Value2 = 1901.5454
$ws.Cells.Item(132, 11).Value2 = 5704.6362
$ws.Cells.Item(132, 13).Value2 = -3174.6362
